$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new drug entry as row 3 of the data source
$ws.Range("A3").Value = "Abreva"

# Page setup touched (adds a pageSetup element, orientation portrait)
$ws.PageSetup.Orientation = 1

# Mirror the saved selection state from the authored workbook
$ws.Range("C4").Select()
